$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): re-order / rename the columns -------------------
# The header labels were reshuffled into a new column order, and a few
# labels were corrected / renamed (DonditionStatus -> ConditionStatus,
# LifeSpan -> LifeSpanYear, GUID -> Guid) as part of the same edit.
$headerRange = $ws.Range("A1:AC1")

$ws.Range("A1").Value = "ProductId"
$ws.Range("B1").Value = "Name"
$ws.Range("C1").Value = "Description"
$ws.Range("D1").Value = "Class"
$ws.Range("E1").Value = "Type"
$ws.Range("F1").Value = "ValidFrom"
$ws.Range("G1").Value = "ValidUntil"
$ws.Range("H1").Value = "GuranteeEnd"
$ws.Range("I1").Value = "LifeSpanYear"
$ws.Range("J1").Value = "ConditionStatus"
$ws.Range("K1").Value = "MaintenanceGroup"
$ws.Range("L1").Value = "MaintenanceRelevance"
$ws.Range("M1").Value = "EquipmentDimension"
$ws.Range("N1").Value = "ProductName"
$ws.Range("O1").Value = "ProductCode"
$ws.Range("P1").Value = "ProductType"
$ws.Range("Q1").Value = "URLLibary"
$ws.Range("R1").Value = "ProducerCode"
$ws.Range("S1").Value = "SupplierCode"
$ws.Range("T1").Value = "Year"
$ws.Range("U1").Value = "SerialNumber"
$ws.Range("V1").Value = "CE-Identification"
$ws.Range("W1").Value = "DeliveryVolume"
$ws.Range("X1").Value = "ContentQuantity"
$ws.Range("Y1").Value = "NetContent"
$ws.Range("Z1").Value = "Payload"
$ws.Range("AA1").Value = "SpaceId"
$ws.Range("AB1").Value = "BusinessPartnerId"
$ws.Range("AC1").Value = "Guid"

# Give the header row its own (distinct) cell format, matching the new
# style slot that was introduced for it.
$headerRange.HorizontalAlignment = 1

# --- Remove the now-unused helper row (row 20) ----------------------------
$ws.Range("A20").Clear()

# --- Sheet view: drop the scrolled-right viewport and select the header --
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
$headerRange.Select()
